$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-262 down to 151-263.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with its data.
$ws.Cells.Item(150, 1).Value2 = 3
$ws.Cells.Item(150, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value2 = "Coquimbo"
$ws.Cells.Item(150, 4).Value2 = 44574
$ws.Cells.Item(150, 5).Value2 = 5
$ws.Cells.Item(150, 6).Value2 = 100112012
$ws.Cells.Item(150, 7).Value2 = "Espinaca"
$ws.Cells.Item(150, 8).Value2 = "Sin especificar"
$ws.Cells.Item(150, 9).Value2 = "Primera"
$ws.Cells.Item(150, 10).Value2 = 155
$ws.Cells.Item(150, 11).Value2 = 4000
$ws.Cells.Item(150, 12).Value2 = 4500
$ws.Cells.Item(150, 13).Value2 = 4274
$ws.Cells.Item(150, 14).Value2 = "$/docena de atados (3 kilos)"
$ws.Cells.Item(150, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(150, 16).Value2 = 1425
$ws.Cells.Item(150, 17).Value2 = 3
$ws.Cells.Item(150, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same number format as the rest of column D.
$ws.Cells.Item(150, 4).NumberFormat = $ws.Cells.Item(151, 4).NumberFormat
